# Generate Report for Handoff
# Adds a new row (row 3) to each of the three sheets (Overview, zh-cn, de-de)
# describing the handoff of the new source file
#   b2777d5e-541d-4972-a24a-a306d7865ae8oo...md
# mirroring the existing row for the 4fb7fe30-... file.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Shared literal values (kept in single-quoted strings so backslashes are
# treated literally and not as escape characters).
# ---------------------------------------------------------------------------
$mdFile        = 'b2777d5e-541d-4972-a24a-a306d7865ae8ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$e2eMdFile     = 'e2e\b2777d5e-541d-4972-a24a-a306d7865ae8ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$status        = 'Ready for handoff'
$handoffDate   = '2016-08-13 00:33:14'
$xlfZhCn       = 'b2777d5e-541d-4972-a24a-a306d7865ae8ooooooooooooooooooooooooooooooooooooooooooo.48183876d04f917072afe446b2539cd2ef30ca8c.zh-cn.xlf'
$handoffZhCn   = '2016-08-13 00:33:07'
$xlfDeDe       = 'b2777d5e-541d-4972-a24a-a306d7865ae8ooooooooooooooooooooooooooooooooooooooooooo.48183876d04f917072afe446b2539cd2ef30ca8c.de-de.xlf'
$handoffDeDe   = $handoffDate

$githubUrl = 'https://github.com/OpenLocalizationTestOrg/oltest/blob/90e9583960f69b0b3d2f51e6a0f4a79ee4a13da8/e2e/' + $mdFile

# ===========================================================================
# Sheet 1: "Overview"
# ===========================================================================
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$rowOverview = $loOverview.ListRows.Add()
$rngOverview = $rowOverview.Range

$rngOverview.Item(1).Value2 = $mdFile        # A3 - File Name
$rngOverview.Item(3).Value2 = ".md"          # C3 - Extension
$rngOverview.Item(4).Value2 = ""             # D3 - Publish URL
$rngOverview.Item(5).Value2 = $status        # E3 - zh-cn
$rngOverview.Item(6).Value2 = $status        # F3 - de-de
$rngOverview.Item(7).Value2 = $handoffDate   # G3 - Latest HO Xliff Generate Date

$hOverview = $wsOverview.Hyperlinks.Add($rngOverview.Item(2), $githubUrl, "", "", $e2eMdFile)

# ===========================================================================
# Sheet 2: "zh-cn"
# ===========================================================================
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$rowZhCn = $loZhCn.ListRows.Add()
$rngZhCn = $rowZhCn.Range

$rngZhCn.Item(2).Value2  = ".md"             # B3 - File Extension
$rngZhCn.Item(3).Value2  = $status           # C3 - Status
$rngZhCn.Item(4).Value2  = "e2e"             # D3 - Source Path
$rngZhCn.Item(5).Value2  = "ht"              # E3 - Priority
$rngZhCn.Item(6).Value2  = "'False"          # F3 - Content Duplicate (force text, not boolean)
$rngZhCn.Item(7).Value2  = $xlfZhCn          # G3 - Latest Handoff File
$rngZhCn.Item(8).Value2  = $handoffZhCn      # H3 - Latest Handoff Datetime
$rngZhCn.Item(9).Value2  = ""                # I3 - Latest Target File
$rngZhCn.Item(10).Value2 = ""                # J3 - Latest Handback File
$rngZhCn.Item(11).Value2 = "0001-01-01 00:00:00"  # K3 - Latest Handback DateTime
$rngZhCn.Item(12).Value2 = ""                # L3 - Reference Tokens
$rngZhCn.Item(13).Value2 = "'True"           # M3 - To be localized (force text, not boolean)
$rngZhCn.Item(14).Value2 = ""                # N3 - Dependency From
$rngZhCn.Item(15).Value2 = "'False"          # O3 - Has metadata (force text, not boolean)
$rngZhCn.Item(16).Value2 = ""                # P3 - Error Detail

$hZhCn = $wsZhCn.Hyperlinks.Add($rngZhCn.Item(1), $githubUrl, "", "", $mdFile)

# ===========================================================================
# Sheet 3: "de-de"
# ===========================================================================
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$rowDeDe = $loDeDe.ListRows.Add()
$rngDeDe = $rowDeDe.Range

$rngDeDe.Item(2).Value2  = ".md"             # B3 - File Extension
$rngDeDe.Item(3).Value2  = $status           # C3 - Status
$rngDeDe.Item(4).Value2  = "e2e"             # D3 - Source Path
$rngDeDe.Item(5).Value2  = "ht"              # E3 - Priority
$rngDeDe.Item(6).Value2  = "'False"          # F3 - Content Duplicate (force text, not boolean)
$rngDeDe.Item(7).Value2  = $xlfDeDe          # G3 - Latest Handoff File
$rngDeDe.Item(8).Value2  = $handoffDeDe      # H3 - Latest Handoff Datetime
$rngDeDe.Item(9).Value2  = ""                # I3 - Latest Target File
$rngDeDe.Item(10).Value2 = ""                # J3 - Latest Handback File
$rngDeDe.Item(11).Value2 = "0001-01-01 00:00:00"  # K3 - Latest Handback DateTime
$rngDeDe.Item(12).Value2 = ""                # L3 - Reference Tokens
$rngDeDe.Item(13).Value2 = "'True"           # M3 - To be localized (force text, not boolean)
$rngDeDe.Item(14).Value2 = ""                # N3 - Dependency From
$rngDeDe.Item(15).Value2 = "'False"          # O3 - Has metadata (force text, not boolean)
$rngDeDe.Item(16).Value2 = ""                # P3 - Error Detail

$hDeDe = $wsDeDe.Hyperlinks.Add($rngDeDe.Item(1), $githubUrl, "", "", $mdFile)

# ---------------------------------------------------------------------------
# Widen the status-ish columns to fit the longer "Ready for handoff" text,
# matching the upstream autofit pass.
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797
$wsZhCn.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsDeDe.Columns.Item(3).ColumnWidth = 17.2159881591797

Write-Output "done"
